# Update cryptos list - Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some D-column values look like plain decimal numbers (single '.') and Excel
# would otherwise auto-convert them from text to a numeric value when set via
# .Value. Force those specific cells to Text format first, then clear the
# formatting afterwards so no stray style index is left behind (these cells
# originally had no explicit style reference at all).
$textForceCells = @("D5","D6","D9","D14","D20","D21","D23","D28","D34","D35","D37","D38","D49")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.982.22"
$ws.Range("E2").Value = "  -2.10%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.416.33"
$ws.Range("E3").Value = "  -1.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "577.92"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6 - Solana
$ws.Range("D6").Value = "152.74"
$ws.Range("E6").Value = "  +3.40%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.27%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "8.02"
$ws.Range("E9").Value = "  +3.51%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.74%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +2.82%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.002.31"
$ws.Range("E12").Value = "  -1.48%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.65%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "28.66"
$ws.Range("E14").Value = "  -2.81%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.423.65"
$ws.Range("E15").Value = "  -1.15%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -0.71%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "62.038.18"
$ws.Range("E17").Value = "  -1.99%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +1.67%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  -0.15%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "8.96"
$ws.Range("E20").Value = "  -4.10%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "381.67"
$ws.Range("E21").Value = "  -1.97%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +0.88%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "75.24"
$ws.Range("E23").Value = "  +0.95%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.03%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.560.06"

# Row 26 - PEPE
$ws.Range("E26").Value = "  -3.96%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -1.99%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "7.67"
$ws.Range("E28").Value = "  +0.00%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.03%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("E30").Value = "  -3.54%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.16%  "

# Row 32 - USDe
$ws.Range("E32").Value = "  -0.04%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -0.81%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "23.22"
$ws.Range("E34").Value = "  -1.07%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +2.97%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -0.58%  "

# Row 37 - Aptos
$ws.Range("D37").Value = "6.95"
$ws.Range("E37").Value = "  -2.93%  "

# Row 38 - Monero
$ws.Range("D38").Value = "168.74"
$ws.Range("E38").Value = "  +0.78%  "

# Row 39 - EnergySwap
$ws.Range("E39").Value = "  -3.35%  "

# Row 40 - RenzoRestakedETH
$ws.Range("D40").Value = "3.450.45"

# Row 41 - Hedera
$ws.Range("E41").Value = "  +1.94%  "

# Row 42 - OKB
$ws.Range("E42").Value = "  +0.67%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -1.79%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +0.57%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  -4.16%  "

# Row 46 - ONDO
$ws.Range("E46").Value = "  -4.36%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.551.11"
$ws.Range("E47").Value = "  -1.53%  "

# Row 48 - Cosmos
$ws.Range("E48").Value = "  +0.19%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "22.58"
$ws.Range("E49").Value = "  -2.08%  "

# Row 50 - dogwifhat
$ws.Range("E50").Value = "  -6.25%  "

# Row 51 - FirstDigitalUSD
$ws.Range("E51").Value = "  +0.07%  "

# Restore default styling on the cells we temporarily forced to Text format,
# so that no stray number-format style index remains on those cells (the
# cells originally had no explicit style reference at all).
foreach ($c in $textForceCells) {
    $ws.Range($c).ClearFormats()
}
